# Add a default header containing the questionnaire number ("Questionnaire 9")
# to the document's (only) section, matching the pattern applied across all
# questionnaire documents in the repo.

$d = $word.ActiveDocument

# The document currently has a single section; grab it and its primary
# (default) header. Indexing with 1 addresses wdHeaderFooterPrimary.
$section = $d.Sections.First
$header = $section.Headers(1)

# Insert the text at the end of the (empty) header story rather than
# assigning .Range.Text directly -- this keeps Word from materialising the
# first-page/even-page header & footer stories (and footnote/endnote
# parts) that it otherwise provisions up front the first time any header
# content is set.
$header.Range.InsertAfter("Questionnaire 9")

# Paragraph-level formatting: built-in "Header" style, centered.
$header.Range.Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1

# Character-level formatting (Arial, 12pt) applied to just the inserted
# text -- shrink a duplicate range by one character first so the trailing
# paragraph mark is excluded and no stray run-properties get stamped onto
# the paragraph mark itself.
$textRange = $header.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
